$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [DateTime]::FromOADate(44279)
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 28000
$ws.Range("L2").Value = 30000
$ws.Range("M2").Value = 29000
$ws.Range("O2").Value = 'Región del Maule'
$ws.Range("P2").Value = 1160
$ws.Range("D3").Value = [DateTime]::FromOADate(44645)
$ws.Range("J3").Value = 220
$ws.Range("K3").Value = 26000
$ws.Range("L3").Value = 27000
$ws.Range("M3").Value = 26455
$ws.Range("P3").Value = 1058
$ws.Range("D4").Value = [DateTime]::FromOADate(44272)
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 22000
$ws.Range("L4").Value = 24000
$ws.Range("M4").Value = 23000
$ws.Range("P4").Value = 920
$ws.Range("D5").Value = [DateTime]::FromOADate(44629)
$ws.Range("J5").Value = 110
$ws.Range("K5").Value = 25000
$ws.Range("L5").Value = 26000
$ws.Range("M5").Value = 25455
$ws.Range("O5").Value = 'Región del Maule'
$ws.Range("P5").Value = 1018
$ws.Range("D6").Value = [DateTime]::FromOADate(44568)
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 26000
$ws.Range("M6").Value = 25500
$ws.Range("O6").Value = 'Región de O''Higgins'
$ws.Range("P6").Value = 1020
$ws.Range("D7").Value = [DateTime]::FromOADate(44320)
$ws.Range("D8").Value = [DateTime]::FromOADate(44587)
$ws.Range("J8").Value = 220
$ws.Range("K8").Value = 23000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 23545
$ws.Range("P8").Value = 942
$ws.Range("D9").Value = [DateTime]::FromOADate(44594)
$ws.Range("J9").Value = 200
$ws.Range("K9").Value = 24000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 24500
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 980
$ws.Range("D11").Value = [DateTime]::FromOADate(44308)
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 28000
$ws.Range("L11").Value = 30000
$ws.Range("M11").Value = 29000
$ws.Range("O11").Value = 'Región del Maule'
$ws.Range("P11").Value = 1160
$ws.Range("D12").Value = [DateTime]::FromOADate(44609)
$ws.Range("J12").Value = 200
$ws.Range("L12").Value = 28000
$ws.Range("M12").Value = 27000
$ws.Range("P12").Value = 1080
$ws.Range("D13").Value = [DateTime]::FromOADate(44194)
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 30000
$ws.Range("L13").Value = 32000
$ws.Range("M13").Value = 31000
$ws.Range("O13").Value = 'Región Metropolitana'
$ws.Range("P13").Value = 1240
$ws.Range("D14").Value = [DateTime]::FromOADate(44574)
$ws.Range("K14").Value = 30000
$ws.Range("L14").Value = 32000
$ws.Range("M14").Value = 31000
$ws.Range("O14").Value = 'Región Metropolitana'
$ws.Range("P14").Value = 1240
$ws.Range("D15").Value = [DateTime]::FromOADate(44244)
$ws.Range("J15").Value = 100
$ws.Range("K15").Value = 25000
$ws.Range("L15").Value = 26000
$ws.Range("M15").Value = 25500
$ws.Range("P15").Value = 1020
$ws.Range("D16").Value = [DateTime]::FromOADate(44552)
$ws.Range("K16").Value = 38000
$ws.Range("L16").Value = 40000
$ws.Range("M16").Value = 39000
$ws.Range("O16").Value = 'Provincia de Limarí'
$ws.Range("P16").Value = 1560
$ws.Range("D17").Value = [DateTime]::FromOADate(44236)
$ws.Range("J17").Value = 100
$ws.Range("M17").Value = 25500
$ws.Range("P17").Value = 1020
$ws.Range("D18").Value = [DateTime]::FromOADate(44216)
$ws.Range("K18").Value = 26000
$ws.Range("L18").Value = 28000
$ws.Range("M18").Value = 27000
$ws.Range("O18").Value = 'Región del Maule'
$ws.Range("P18").Value = 1080
$ws.Range("D19").Value = [DateTime]::FromOADate(44210)
$ws.Range("K19").Value = 32000
$ws.Range("L19").Value = 34000
$ws.Range("M19").Value = 33000
$ws.Range("P19").Value = 1320
$ws.Range("D20").Value = [DateTime]::FromOADate(44602)
$ws.Range("J20").Value = 130
$ws.Range("K20").Value = 20000
$ws.Range("L20").Value = 21000
$ws.Range("M20").Value = 20385
$ws.Range("P20").Value = 815
$ws.Range("D21").Value = [DateTime]::FromOADate(44636)
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 22000
$ws.Range("L21").Value = 23000
$ws.Range("M21").Value = 22375
$ws.Range("P21").Value = 895
$ws.Range("D22").Value = [DateTime]::FromOADate(44203)
$ws.Range("K22").Value = 25000
$ws.Range("L22").Value = 26000
$ws.Range("M22").Value = 25500
$ws.Range("O22").Value = 'Región de O''Higgins'
$ws.Range("P22").Value = 1020
$ws.Range("D23").Value = [DateTime]::FromOADate(44651)
$ws.Range("J23").Value = 250
$ws.Range("K23").Value = 28000
$ws.Range("L23").Value = 30000
$ws.Range("M23").Value = 28960
$ws.Range("O23").Value = 'Región Metropolitana'
$ws.Range("P23").Value = 1158
$ws.Range("D24").Value = [DateTime]::FromOADate(44603)
$ws.Range("J24").Value = 130
$ws.Range("K24").Value = 25000
$ws.Range("L24").Value = 26000
$ws.Range("M24").Value = 25385
$ws.Range("O24").Value = 'Región del Maule'
$ws.Range("P24").Value = 1015
$ws.Range("D25").Value = [DateTime]::FromOADate(44316)
$ws.Range("K25").Value = 26000
$ws.Range("L25").Value = 27000
$ws.Range("M25").Value = 26500
$ws.Range("P25").Value = 1060
$ws.Range("D26").Value = [DateTime]::FromOADate(44631)
$ws.Range("J26").Value = 110
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 21000
$ws.Range("M26").Value = 20455
$ws.Range("O26").Value = 'Región de O''Higgins'
$ws.Range("P26").Value = 818
$ws.Range("D27").Value = [DateTime]::FromOADate(44642)
$ws.Range("J27").Value = 120
$ws.Range("K27").Value = 20000
$ws.Range("L27").Value = 20000
$ws.Range("M27").Value = 20000
$ws.Range("P27").Value = 800
$ws.Range("D28").Value = [DateTime]::FromOADate(44560)
$ws.Range("K28").Value = 25000
$ws.Range("L28").Value = 26000
$ws.Range("M28").Value = 25500
$ws.Range("P28").Value = 1020
$ws.Range("D29").Value = [DateTime]::FromOADate(44342)
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 28000
$ws.Range("L29").Value = 30000
$ws.Range("M29").Value = 29000
$ws.Range("O29").Value = 'Región Metropolitana'
$ws.Range("P29").Value = 1160
$ws.Range("D30").Value = [DateTime]::FromOADate(44251)
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 27000
$ws.Range("M30").Value = 27500
$ws.Range("P30").Value = 1100
$ws.Range("D31").Value = [DateTime]::FromOADate(44230)
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 35000
$ws.Range("L31").Value = 36000
$ws.Range("M31").Value = 35500
$ws.Range("O31").Value = 'Región del Maule'
$ws.Range("P31").Value = 1420
$ws.Range("D32").Value = [DateTime]::FromOADate(44328)
$ws.Range("J32").Value = 100
$ws.Range("K32").Value = 32000
$ws.Range("L32").Value = 34000
$ws.Range("M32").Value = 33000
$ws.Range("O32").Value = 'Región Metropolitana'
$ws.Range("P32").Value = 1320
$ws.Range("D33").Value = [DateTime]::FromOADate(44644)
$ws.Range("J33").Value = 130
$ws.Range("L33").Value = 21000
$ws.Range("M33").Value = 20615
$ws.Range("P33").Value = 825
$ws.Range("D34").Value = [DateTime]::FromOADate(44294)
$ws.Range("J34").Value = 100
$ws.Range("L34").Value = 22000
$ws.Range("M34").Value = 21000
$ws.Range("O34").Value = 'Región del Maule'
$ws.Range("P34").Value = 840
$ws.Range("D35").Value = [DateTime]::FromOADate(44265)
$ws.Range("K35").Value = 22000
$ws.Range("L35").Value = 24000
$ws.Range("M35").Value = 23000
$ws.Range("P35").Value = 920
$ws.Range("D36").Value = [DateTime]::FromOADate(44624)
$ws.Range("J36").Value = 150
$ws.Range("K36").Value = 25000
$ws.Range("L36").Value = 26000
$ws.Range("M36").Value = 25467
$ws.Range("O36").Value = 'Región Metropolitana'
$ws.Range("P36").Value = 1019
$ws.Range("D37").Value = [DateTime]::FromOADate(44313)
$ws.Range("K37").Value = 30000
$ws.Range("L37").Value = 32000
$ws.Range("M37").Value = 31000
$ws.Range("O37").Value = 'Región Metropolitana'
$ws.Range("P37").Value = 1240
$ws.Range("D38").Value = [DateTime]::FromOADate(44558)
$ws.Range("J38").Value = 250
$ws.Range("K38").Value = 15000
$ws.Range("L38").Value = 16000
$ws.Range("M38").Value = 15400
$ws.Range("O38").Value = 'Región Metropolitana'
$ws.Range("P38").Value = 616
$ws.Range("D41").Value = [DateTime]::FromOADate(44638)
$ws.Range("J41").Value = 220
$ws.Range("K41").Value = 20000
$ws.Range("L41").Value = 22000
$ws.Range("M41").Value = 21091
$ws.Range("O41").Value = 'Región del Maule'
$ws.Range("P41").Value = 844
